# Fruta / hortaliza, semanal
# Insert two new weekly observation rows at the top of the data block
# (row 166 onward), pushing all existing rows down by two, then populate
# the two new rows with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 166..232 down to 168..234, leaving two blank rows at 166-167.
$ws.Rows.Item(166).Resize(2).Insert()

# New row 166: Especial quality observation
$ws.Cells.Item(166, 1).Value2 = 7
$ws.Cells.Item(166, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value2 = "Ñuble"
$ws.Cells.Item(166, 4).Value2 = 44825
$ws.Cells.Item(166, 5).Value2 = 16
$ws.Cells.Item(166, 6).Value2 = "Fruta"
$ws.Cells.Item(166, 7).Value2 = 100104
$ws.Cells.Item(166, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(166, 9).Value2 = 100104005
$ws.Cells.Item(166, 10).Value2 = "Pera"
$ws.Cells.Item(166, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(166, 12).Value2 = "Especial"
$ws.Cells.Item(166, 13).Value2 = 40
$ws.Cells.Item(166, 14).Value2 = 12000
$ws.Cells.Item(166, 15).Value2 = 12000
$ws.Cells.Item(166, 16).Value2 = 12000
$ws.Cells.Item(166, 17).Value2 = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(166, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(166, 19).Value2 = 750
$ws.Cells.Item(166, 20).Value2 = 16

# New row 167: Primera quality observation
$ws.Cells.Item(167, 1).Value2 = 7
$ws.Cells.Item(167, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(167, 3).Value2 = "Ñuble"
$ws.Cells.Item(167, 4).Value2 = 44825
$ws.Cells.Item(167, 5).Value2 = 16
$ws.Cells.Item(167, 6).Value2 = "Fruta"
$ws.Cells.Item(167, 7).Value2 = 100104
$ws.Cells.Item(167, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(167, 9).Value2 = 100104005
$ws.Cells.Item(167, 10).Value2 = "Pera"
$ws.Cells.Item(167, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(167, 12).Value2 = "Primera"
$ws.Cells.Item(167, 13).Value2 = 80
$ws.Cells.Item(167, 14).Value2 = 10500
$ws.Cells.Item(167, 15).Value2 = 11000
$ws.Cells.Item(167, 16).Value2 = 10750
$ws.Cells.Item(167, 17).Value2 = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(167, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(167, 19).Value2 = 672
$ws.Cells.Item(167, 20).Value2 = 16
